# Decrement the "剩余" (remaining) value in column E for each data row
# (rows 2-99) by 1, except row 36 which stays unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E is the 5th column
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value = $current - 1
    }
}
